$d = $word.ActiveDocument

# 1. Description paragraph: "any human user." -> "any medical professional or app user."
$d.Content.Find.Execute(
    "any human user.", $true, $false, $false, $false, $false,
    $true, 1, $false, "any medical professional or app user.", 2) | Out-Null

# 2. Actors paragraph: "and the user of the system." -> "and the general public user of the system."
$d.Content.Find.Execute(
    "and the user of the system.", $true, $false, $false, $false, $false,
    $true, 1, $false, "and the general public user of the system.", 2) | Out-Null

# 3. Pre-conditions paragraph: remove stray " result" word
$d.Content.Find.Execute(
    "available result: hospital", $true, $false, $false, $false, $false,
    $true, 1, $false, "available: hospital", 2) | Out-Null

# 4. Post-conditions paragraph: amber -> red, and "coming in" -> "entering ... central"
$d.Content.Find.Execute(
    "and get a new amber status. This notification will be sent and received by the phones only in order to avoid peaks of information coming in the track and trace database.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "and get a new red status. This notification will be sent and received by the phones only in order to avoid peaks of information entering the track and trace central database.",
    2) | Out-Null

# 5. Relocate the _GoBack bookmark from the "Alternative Scenarios" heading
#    to just after "database" (before the final ". ") in the Post-conditions paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$rngAfterDatabase = $d.Content
$rngAfterDatabase.Find.Execute(
    "central database", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$bookmarkRange = $d.Range($rngAfterDatabase.End, $rngAfterDatabase.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# 6. Main Scenario step 3: "The medical professionals update" -> "The medical system updates"
$d.Content.Find.Execute(
    "The medical professionals update the track and trace system",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The medical system updates the track and trace system",
    2) | Out-Null
